# Insert a new weekly price record for "Ajo / Chino / Primera" at row 128.
# This pushes every existing record at/after row 128 down by one row
# (old row 128 -> new row 129, ..., old row 226 -> new row 227), and the
# sheet's used range grows from A1:R226 to A1:R227.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 128:226 down to 129:227 by inserting a blank row at 128.
$ws.Rows("128:128").Insert()

# Populate the newly inserted row 128 with the new record.
$ws.Cells.Item(128, 1).Value  = 7
$ws.Cells.Item(128, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(128, 3).Value  = "Ñuble"
$ws.Cells.Item(128, 4).Value  = 44673
$ws.Cells.Item(128, 5).Value  = 16
$ws.Cells.Item(128, 6).Value  = 100112003
$ws.Cells.Item(128, 7).Value  = "Ajo"
$ws.Cells.Item(128, 8).Value  = "Chino"
$ws.Cells.Item(128, 9).Value  = "Primera"
$ws.Cells.Item(128, 10).Value = 30
$ws.Cells.Item(128, 11).Value = 21000
$ws.Cells.Item(128, 12).Value = 22000
$ws.Cells.Item(128, 13).Value = 21500
$ws.Cells.Item(128, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(128, 15).Value = "China"
$ws.Cells.Item(128, 16).Value = 2150
$ws.Cells.Item(128, 17).Value = 10
$ws.Cells.Item(128, 18).Value = "Hortaliza"
